$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = 4.7
$ws.Range("I2").Value = 4.9
$ws.Range("J2").Value = 4.5
$ws.Range("K2").Value = 4.6
$ws.Range("P2").Value = 3.15
$ws.Range("R2").Value = 1.88
$ws.Range("S2").Value = 2.08
$ws.Range("T2").Value = 1.51
$ws.Range("U2").Value = 2.84
$ws.Range("X2").Value = 36
$ws.Range("Z2").Value = 46
$ws.Range("AD2").Value = 19.5
$ws.Range("AI2").Value = 40
$ws.Range("AN2").Value = 5.9
$ws.Range("F3").Value = 1.93
$ws.Range("H3").Value = 4.2
$ws.Range("M3").Value = 1.05
$ws.Range("U3").Value = 2.46
$ws.Range("X3").Value = 19.5
$ws.Range("AI3").Value = 46
$ws.Range("G4").Value = 1.94
$ws.Range("J4").Value = 4
$ws.Range("K4").Value = 4.7
$ws.Range("N4").Value = 4.8
$ws.Range("P4").Value = 2.34
$ws.Range("Q4").Value = 1.6
$ws.Range("W4").Value = 2.06
$ws.Range("Y4").Value = 25
$ws.Range("J5").Value = 2.9
$ws.Range("K5").Value = 3.55
$ws.Range("L5").Value = 1.35
$ws.Range("P5").Value = 1.62
$ws.Range("S5").Value = 3.5
$ws.Range("F6").Value = 1.49
$ws.Range("G6").Value = 1.83
$ws.Range("I6").Value = 22
$ws.Range("J6").Value = 3.7
$ws.Range("K6").Value = 7.2
$ws.Range("L6").Value = 1.34
$ws.Range("M6").Value = 1.06
$ws.Range("P6").Value = 1.76
$ws.Range("R6").Value = 1.23
$ws.Range("S6").Value = 3
$ws.Range("V6").Value = 1.06
$ws.Range("W6").Value = 2.18
$ws.Range("F8").Value = 2.14
$ws.Range("H8").Value = 1.01
$ws.Range("K8").Value = 4
$ws.Range("Q8").Value = 1.6
$ws.Range("S8").Value = 1.6
$ws.Range("V8").Value = 1.27
$ws.Range("F9").Value = 2.28
$ws.Range("G9").Value = 2.3
$ws.Range("H9").Value = 3.6
$ws.Range("I9").Value = 3.65
$ws.Range("O9").Value = 1.38
$ws.Range("V9").Value = 1.37
$ws.Range("W9").Value = 1.76
$ws.Range("AD9").Value = 14.5
$ws.Range("AE9").Value = 44
$ws.Range("AJ9").Value = 29
$ws.Range("AN9").Value = 21
$ws.Range("S10").Value = 2.68
$ws.Range("T10").Value = 1.78
$ws.Range("AD10").Value = 24
$ws.Range("AE10").Value = 80
$ws.Range("G11").Value = 2.06
$ws.Range("H11").Value = 3.6
$ws.Range("J11").Value = 4
$ws.Range("W11").Value = 1.94
$ws.Range("F12").Value = 2.28
$ws.Range("G12").Value = 2.3
$ws.Range("H12").Value = 3.3
$ws.Range("I12").Value = 3.35
$ws.Range("J12").Value = 3.8
$ws.Range("K12").Value = 3.85
$ws.Range("N12").Value = 4.1
$ws.Range("T12").Value = 1.78
$ws.Range("U12").Value = 2.22
$ws.Range("V12").Value = 1.42
$ws.Range("W12").Value = 1.76
$ws.Range("AE12").Value = 36
$ws.Range("AF12").Value = 13.5
$ws.Range("AI12").Value = 46
$ws.Range("AM12").Value = 85
$ws.Range("AN12").Value = 17.5
$ws.Range("H13").Value = 1.41
$ws.Range("I13").Value = 1.42
$ws.Range("N13").Value = 5.2
$ws.Range("T13").Value = 1.95
$ws.Range("Z13").Value = 8.2